$wb = $excel.ActiveWorkbook

# Map of row -> new "want to go" count, applies to the 展览 sheet (rows 2-18)
# and to the 全部类型 sheet (same rows 2-18, since it mirrors 展览's rows).
$countUpdates = @{
    3  = 275
    5  = 831
    7  = 6798
    8  = 58
    10 = 121
    11 = 87
    13 = 13
    15 = 21
    16 = 231
    17 = 586
    18 = 67
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $eCell = $ws.Cells.Item($r, 5)
        $eText = $eCell.Value()
        if ($eText -ne $null -and $eText -ne "") {
            $newText = $eText -replace "-", " - "
            if ($newText -ne $eText) {
                $eCell.Value = $newText
            }
        }

        if ($countUpdates.ContainsKey($r)) {
            $fCell = $ws.Cells.Item($r, 6)
            $fVal = $fCell.Value()
            if ($fVal -ne $null -and $fVal -ne "") {
                $fCell.Value = $countUpdates[$r]
            }
        }
    }
}
